# feat: add 2022-Q1 data
#
# 1. Insert a new worksheet "2022-Q1" between "2021-Q4" and "总计".
# 2. Populate "2022-Q1" with the per-fund holdings for that quarter
#    (mirrors the layout already used by the "2021-Q4" sheet).
# 3. Update the "总计" (totals) sheet: push the existing 2021-Q4 summary
#    row down one row and add a new top row summarizing 2022-Q1.

$wb = $excel.ActiveWorkbook

$q4Sheet = $wb.Worksheets.Item("2021-Q4")

# --- 1) Create the new sheet right after "2021-Q4" ------------------------
$q1Sheet = $wb.Worksheets.Add($null, $q4Sheet)
$q1Sheet.Name = "2022-Q1"

# NOTE: Worksheets.Item(...) resolves by current tab position, not a
# stable sheet identity, so fetch "总计" only *after* the insert above has
# shifted it from index 2 to index 3 - otherwise this handle would silently
# alias the newly-inserted "2022-Q1" sheet instead.
$totalSheet = $wb.Worksheets.Item("总计")

# --- 2) Seed "2022-Q1" from "2021-Q4" so it inherits identical cell
#        styles (bold+boxed header B1:H1, boxed index column A2:A3), then
#        overwrite the actual values/text below.
$q4Sheet.Range("B1:H3").Copy($q1Sheet.Range("B1"))
$q4Sheet.Range("A2:A3").Copy($q1Sheet.Range("A2"))

# Headers (row 1)
$q1Sheet.Range("B1").Value = "基金代码"
$q1Sheet.Range("C1").Value = "基金名称"
$q1Sheet.Range("D1").Value = "基金规模"
$q1Sheet.Range("E1").Value = "股票总仓位"
$q1Sheet.Range("F1").Value = "仓位占比"
$q1Sheet.Range("G1").Value = "持有市值(亿元)"
$q1Sheet.Range("H1").Value = "仓位排名"

# Columns B, D, E, F, G hold numeric-looking text (fund codes with leading
# zeros, percentages, etc.) that must stay plain text rather than being
# coerced to numbers - force text format, write, then drop the format again
# so the cells end up as plain unstyled text cells (matching "2021-Q4").
$textRangeB = $q1Sheet.Range("B2:B3")
$textRangeDG = $q1Sheet.Range("D2:G3")
$textRangeB.NumberFormat = "@"
$textRangeDG.NumberFormat = "@"

$q1Sheet.Range("A2").Value = 0
$q1Sheet.Range("B2").Value = "001917"
$q1Sheet.Range("C2").Value = "招商量化精选股票A"
$q1Sheet.Range("D2").Value = "2.33"
$q1Sheet.Range("E2").Value = "94.20"
$q1Sheet.Range("F2").Value = "1.52"
$q1Sheet.Range("G2").Value = "0.0354"
$q1Sheet.Range("H2").Value = 2

$q1Sheet.Range("A3").Value = 1
$q1Sheet.Range("B3").Value = "007950"
$q1Sheet.Range("C3").Value = "招商量化精选股票C"
$q1Sheet.Range("D3").Value = "0.56"
$q1Sheet.Range("E3").Value = "94.20"
$q1Sheet.Range("F3").Value = "1.52"
$q1Sheet.Range("G3").Value = "0.0085"
$q1Sheet.Range("H3").Value = 2

$textRangeB.ClearFormats()
$textRangeDG.ClearFormats()

# --- 3) Update "总计": shift the existing 2021-Q4 row down to row 3, then
#        write the new 2022-Q1 summary into row 2. -------------------------
$totalSheet.Range("A2").Copy($totalSheet.Range("A3"))
$totalSheet.Range("A3").Value = 1
$totalSheet.Range("B3").Value = "2021-Q4"
$totalSheet.Range("C3").Value = 4
$totalSheet.Range("D3").Value = 0.06

$totalSheet.Range("A2").Value = 0
$totalSheet.Range("B2").Value = "2022-Q1"
$totalSheet.Range("C2").Value = 2
$totalSheet.Range("D2").Value = 0.04

# Restore the originally active tab (the new-sheet insert/rename above left
# "2022-Q1" selected).
$q4Sheet.Activate()
